# Add two new batches (30 and 31) of device-master test rows, mirroring
# the existing data pattern (Finger Print Scanner / IRIS Scanner / Web
# Camera / Document Scanner / Printer, each with its own MAC address and
# serial number), appended as worksheet rows 147-156.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# id, name, mac_address, serial_num, dspec_id
$newRows = @(
    ,@(3000166, "Finger Print Scanner 30", "D6-15-AC-80-6B-86", "BS563Q2230814", 165)
    ,@(3000167, "IRIS Scanner 30",         "6D-58-E2-DF-74-34", "BS563Q2230815", 327)
    ,@(3000168, "Web Camera 30",           "E2-A8-56-86-15-30", "BS563Q2230816", 736)
    ,@(3000169, "Document Scanner 30",     "72-E8-B9-FD-63-65", "BS563Q2230817", 801)
    ,@(3000170, "Printer 30",              "D3-F3-A4-50-AD-12", "BS563Q2230818", 920)
    ,@(3000171, "Finger Print Scanner 31", "06-16-D0-0B-A6-E4", "BS563Q2230819", 165)
    ,@(3000172, "IRIS Scanner 31",         "21-78-45-AC-E9-20", "BS563Q2230820", 327)
    ,@(3000173, "Web Camera 31",           "3C-E8-87-99-DB-FA", "BS563Q2230821", 736)
    ,@(3000174, "Document Scanner 31",     "BF-55-53-98-40-08", "BS563Q2230822", 801)
    ,@(3000175, "Printer 31",              "5A-43-36-46-22-EB", "BS563Q2230823", 920)
)

$r = 147
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]          # A - id
    $ws.Cells.Item($r, 2).Value = $row[1]          # B - name
    $ws.Cells.Item($r, 3).Value = $row[2]          # C - mac_address
    $ws.Cells.Item($r, 4).Value = $row[3]          # D - serial_num
    # column E (ip_address) is left blank, matching the existing rows
    $ws.Cells.Item($r, 6).Value = $row[4]          # F - dspec_id
    $ws.Cells.Item($r, 7).Value = "eng"            # G - lang_code
    $ws.Cells.Item($r, 8).Value = $true            # H - is_active
    $ws.Cells.Item($r, 8).HorizontalAlignment = -4131  # xlLeft, matches style of existing H cells
    $ws.Cells.Item($r, 9).Value = "superadmin"     # I - cr_by
    $ws.Cells.Item($r, 10).Value = "now()"         # J - cr_dtimes
    $r++
}

# Move the selection to reflect where editing left off, as Excel would.
$ws.Range("E156").Select() | Out-Null
